# Schema_colors.xlsx - "Towards code cleaning and heating/cooling calculation"
#
# 1) Update the Coal color swatch (row 7, column B) from the old near-black
#    (#231f20) to the new dark grey (#35373a).
# 2) Update the Renewable color swatch (row 24, column B) from the teal
#    (#00af8c, still used elsewhere for Kea/renewable rows) to the new
#    lighter mint (#8dd2c0), so "Renewable" gets its own distinct color.
# 3) Move the active selection from F15 to E20 and scroll the view down so
#    row 10 is at the top of the visible window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Coal color: #231f20 -> #35373a -------------------------------------
$ws.Range("B7").Value = "#35373a"

# --- Renewable color: #00af8c -> #8dd2c0 ---------------------------------
$ws.Range("B24").Value = "#8dd2c0"

# --- Update the view: scroll so row 10 is the top row, select E20 -------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("E20").Select() | Out-Null
